# Add DateTimeTest test cases: a row for "empty cells", a row for
# "explicit null cells" (with literal "null" strings in B:E), and a new
# column F filled with "test" as a marker/sentinel column down through
# row 4. Finally move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: empty-cells test case (only column A populated)
$ws.Range("A3").Value = "empty cells"

# Row 4: explicit-null-cells test case
$ws.Range("A4").Value = "explicit null cells"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"

# New column F, marker value "test" for header + all data rows
$ws.Range("F1").Value = "test"
$ws.Range("F2").Value = "test"
$ws.Range("F3").Value = "test"
$ws.Range("F4").Value = "test"

# Match the saved selection from the authored workbook
$ws.Range("A2").Select()
